$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text does not look like a plain number (e.g. the price
# column uses "." as a thousands separator, so values like "55.735.19" or
# "2.341.32" read back fine as text on their own, and the "Volume(1h)"
# column always has a trailing '%' plus padding spaces). These can be
# assigned directly.
$textSafeUpdates = @{
    "D2" = '55.735.19'
    "E2" = '  -2.24%  '
    "D3" = '2.341.32'
    "E3" = '  -2.44%  '
    "E4" = '  -0.15%  '
    "E5" = '  -0.60%  '
    "E6" = '  -3.72%  '
    "E7" = '  -0.03%  '
    "E8" = '  -3.41%  '
    "D9" = '2.348.43'
    "E9" = '  -2.49%  '
    "E10" = '  -0.24%  '
    "E11" = '  -0.33%  '
    "E12" = '  +3.21%  '
    "E13" = '  -1.86%  '
    "D14" = '2.756.20'
    "E14" = '  -2.52%  '
    "D15" = '55.692.05'
    "E15" = '  -2.21%  '
    "E16" = '  -1.12%  '
    "E17" = '  -2.49%  '
    "D18" = '2.323.15'
    "E18" = '  -2.47%  '
    "E19" = '  -3.27%  '
    "E20" = '  -0.38%  '
    "E21" = '  -2.39%  '
    "E22" = '  -1.07%  '
    "E23" = '  -0.06%  '
    "E24" = '  -3.73%  '
    "E25" = '  +0.10%  '
    "E26" = '  -1.32%  '
    "E27" = '  -3.94%  '
    "E28" = '  -4.79%  '
    "E29" = '  -2.46%  '
    "E30" = '  -2.01%  '
    "D31" = '0.0₃0703'
    "E31" = '  -3.58%  '
    "E32" = '  -0.02%  '
    "E33" = '  -1.66%  '
    "E35" = '  -5.39%  '
    "E36" = '  -2.04%  '
    "E37" = '  -2.68%  '
    "E38" = '  -5.17%  '
    "E39" = '  -1.60%  '
    "E40" = '  -2.15%  '
    "E41" = '  -4.75%  '
    "E42" = '  -1.33%  '
    "E43" = '  -4.18%  '
    "E44" = '  -5.58%  '
    "E45" = '  -2.99%  '
    "E46" = '  -2.43%  '
    "E47" = '  -6.33%  '
    "E48" = '  -2.75%  '
    "E49" = '  -2.77%  '
    "E50" = '  -2.13%  '
    "E51" = '  +0.09%  '
}

# Cells whose new text WOULD be silently re-interpreted by Excel as a
# real number (e.g. "4.79", "0.999", "21.60"), which would lose the
# exact original text (trailing zeros, etc.). These columns store text
# in the workbook, so force Text format on the cell before assigning
# to keep the value as a string.
$textForcedUpdates = @{
    "D5" = '503.97'
    "D6" = '128.59'
    "D12" = '4.79'
    "D13" = '0.319'
    "D16" = '21.60'
    "D20" = '309.90'
    "D21" = '3.98'
    "D22" = '6.20'
    "D23" = '0.999'
    "D24" = '65.32'
    "D28" = '7.08'
    "D29" = '171.83'
    "D30" = '1.64'
    "D33" = '5.78'
    "D36" = '17.62'
    "D38" = '3.64'
    "D39" = '0.817'
    "D43" = '126.35'
    "D44" = '4.70'
    "D47" = '236.42'
    "D48" = '0.0476'
    "D50" = '16.81'
    "D51" = '0.953'
}

foreach ($cellRef in $textSafeUpdates.Keys) {
    $ws.Range($cellRef).Value = $textSafeUpdates[$cellRef]
}

foreach ($cellRef in $textForcedUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$cellRef]
}
